# Atualizando planilha modelo para importacao de militares.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the MATRICULA column (B2:B26) from the old random "REGxxx"
# codes to sequential "REG00001".."REG00025" codes.
for ($i = 2; $i -le 26; $i++) {
    $n = $i - 1
    $code = "REG" + $n.ToString().PadLeft(5, '0')
    $ws.Cells.Item($i, 2).Value = $code
}

# Move the active selection from A1 to C28.
$ws.Range("C28").Select() | Out-Null
